$wb = $excel.ActiveWorkbook
$rng = $wb.Worksheets.Item("Paises").Range("A2:A11")
$nm = "África"
Write-Output "Length: $($nm.Length)"
$rng.Name = $nm
